$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 (shifts existing data rows down by one)
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).ClearFormats()

# Match the date-format style used by the other rows' Fecha column
$ws.Range("D2").NumberFormat = $ws.Range("D3").NumberFormat

# Fill in the new week's record
$ws.Range("A2").Value = 10
$ws.Range("B2").Value = "Vega Modelo de Temuco"
$ws.Range("C2").Value = "La Araucanía"
$ws.Range("D2").Value = 44483
$ws.Range("E2").Value = 9
$ws.Range("F2").Value = 100112042
$ws.Range("G2").Value = "Locoto"
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 50
$ws.Range("K2").Value = 2200
$ws.Range("L2").Value = 2200
$ws.Range("M2").Value = 2200
$ws.Range("N2").Value = "$/kilo"
$ws.Range("O2").Value = "Región de Arica y Parinacota"
$ws.Range("P2").Value = 2200
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = "Hortaliza"
